$wb = $excel.ActiveWorkbook

# Update data: IAM sheet E4/E5 "SYN_prod" -> "SYN_nonprod"
$iam = $wb.Worksheets.Item("IAM")
$iam.Range("E4").Value = "SYN_nonprod"
$iam.Range("E5").Value = "SYN_nonprod"

# Select E5 on IAM sheet, and make IAM the active sheet
$iam.Range("E5").Select()
$iam.Activate()
